$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = "Jeremi.Beld@live.com"
$ws.Range("B47").Value = 62267877
